# Scheduled runner update: refresh market-board derived values
# (currentAveragePrice*, LevePrice*, LeveProfit* columns) across all
# Disciple of the Hand job sheets in the Adamantoise profits workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J17").Value = 2012038.6
$ws.Range("L17").Value = 6036115.800000001
$ws.Range("N17").Value = -6036451.800000001
$ws.Range("H51").Value = 3743
$ws.Range("J51").Value = 2940.3333
$ws.Range("L51").Value = 2940.3333
$ws.Range("N51").Value = -3908.3333
$ws.Range("H106").Value = 11112104
$ws.Range("I106").Value = 11112104
$ws.Range("K106").Value = 11112104
$ws.Range("M106").Value = -11111473
$ws.Range("H129").Value = 2199.7273
$ws.Range("I129").Value = 998.5
$ws.Range("K129").Value = 2995.5
$ws.Range("M129").Value = 2004.5
$ws.Range("H132").Value = 5351.276
$ws.Range("I132").Value = 5414.3335
$ws.Range("K132").Value = 16243.0005
$ws.Range("M132").Value = -13713.0005
$ws.Range("H137").Value = 33305.777
$ws.Range("I137").Value = 45109.58
$ws.Range("K137").Value = 135328.74
$ws.Range("M137").Value = -132778.74
$ws.Range("H138").Value = 3009.634
$ws.Range("J138").Value = 3673.1538
$ws.Range("L138").Value = 11019.4614
$ws.Range("N138").Value = -21299.4614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5721.7354
$ws.Range("I32").Value = 1921.6
$ws.Range("J32").Value = 16277.667
$ws.Range("K32").Value = 1921.6
$ws.Range("L32").Value = 16277.667
$ws.Range("M32").Value = -1634.6
$ws.Range("N32").Value = -16851.667
$ws.Range("H61").Value = 3543.4546
$ws.Range("I61").Value = 3347.8
$ws.Range("K61").Value = 3347.8
$ws.Range("M61").Value = -3135.8
$ws.Range("H74").Value = 2539.8696
$ws.Range("I74").Value = 2234.1428
$ws.Range("K74").Value = 2234.1428
$ws.Range("M74").Value = -1360.1428
$ws.Range("H77").Value = 2539.8696
$ws.Range("I77").Value = 2234.1428
$ws.Range("K77").Value = 11170.714
$ws.Range("M77").Value = -6802.714
$ws.Range("H102").Value = 2641.9333
$ws.Range("I102").Value = 2363
$ws.Range("K102").Value = 2363
$ws.Range("M102").Value = -741
$ws.Range("H110").Value = 1550.7916
$ws.Range("I110").Value = 1145.5625
$ws.Range("J110").Value = 2361.25
$ws.Range("K110").Value = 1145.5625
$ws.Range("L110").Value = 2361.25
$ws.Range("M110").Value = 899.4375
$ws.Range("N110").Value = -6451.25
$ws.Range("H122").Value = 4934.8374
$ws.Range("I122").Value = 2142.625
$ws.Range("J122").Value = 8461.842000000001
$ws.Range("K122").Value = 6427.875
$ws.Range("L122").Value = 25385.526
$ws.Range("M122").Value = -3977.875
$ws.Range("N122").Value = -30285.526
$ws.Range("H132").Value = 457594.53
$ws.Range("I132").Value = 529107.25
$ws.Range("J132").Value = 4680.6665
$ws.Range("K132").Value = 1587321.75
$ws.Range("L132").Value = 14041.9995
$ws.Range("M132").Value = -1584791.75
$ws.Range("N132").Value = -19101.9995
$ws.Range("H136").Value = 3543.4546
$ws.Range("I136").Value = 3347.8
$ws.Range("K136").Value = 10043.4
$ws.Range("M136").Value = -7493.400000000001
$ws.Range("H141").Value = 559999.5
$ws.Range("J141").Value = 559999.5
$ws.Range("L141").Value = 559999.5
$ws.Range("N141").Value = -570359.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 11907116
$ws.Range("I134").Value = 1859.3182
$ws.Range("J134").Value = 55559724
$ws.Range("K134").Value = 5577.9546
$ws.Range("L134").Value = 166679172
$ws.Range("M134").Value = -3042.9546
$ws.Range("N134").Value = -166684242

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3417.743
$ws.Range("I31").Value = 2516
$ws.Range("J31").Value = 4094.05
$ws.Range("K31").Value = 2516
$ws.Range("L31").Value = 4094.05
$ws.Range("M31").Value = -2221
$ws.Range("N31").Value = -4684.05
$ws.Range("H34").Value = 3417.743
$ws.Range("I34").Value = 2516
$ws.Range("J34").Value = 4094.05
$ws.Range("K34").Value = 2516
$ws.Range("L34").Value = 4094.05
$ws.Range("M34").Value = -2314
$ws.Range("N34").Value = -4498.05
$ws.Range("H58").Value = 3375.4
$ws.Range("I58").Value = 3212.25
$ws.Range("J58").Value = 3665.4443
$ws.Range("K58").Value = 3212.25
$ws.Range("L58").Value = 3665.4443
$ws.Range("M58").Value = -3009.25
$ws.Range("N58").Value = -4071.4443
$ws.Range("H105").Value = 3269.8572
$ws.Range("J105").Value = 3872.25
$ws.Range("L105").Value = 3872.25
$ws.Range("N105").Value = -7366.25
$ws.Range("H132").Value = 3319.8076
$ws.Range("I132").Value = 2832.6365
$ws.Range("J132").Value = 5999.25
$ws.Range("K132").Value = 8497.9095
$ws.Range("L132").Value = 17997.75
$ws.Range("M132").Value = -5967.9095
$ws.Range("N132").Value = -23057.75
$ws.Range("H134").Value = 2418.5557
$ws.Range("I134").Value = 1726.069
$ws.Range("K134").Value = 5178.207
$ws.Range("M134").Value = -2643.207
$ws.Range("H136").Value = 3375.4
$ws.Range("I136").Value = 3212.25
$ws.Range("J136").Value = 3665.4443
$ws.Range("K136").Value = 9636.75
$ws.Range("L136").Value = 10996.3329
$ws.Range("M136").Value = -7086.75
$ws.Range("N136").Value = -16096.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 2437.5
$ws.Range("J88").Value = 2437.5
$ws.Range("L88").Value = 7312.5
$ws.Range("N88").Value = -8168.5
$ws.Range("H91").Value = 2437.5
$ws.Range("J91").Value = 2437.5
$ws.Range("L91").Value = 7312.5
$ws.Range("N91").Value = -10276.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4399.6
$ws.Range("I70").Value = 3875
$ws.Range("K70").Value = 3875
$ws.Range("M70").Value = -3605
$ws.Range("H73").Value = 4399.6
$ws.Range("I73").Value = 3875
$ws.Range("K73").Value = 3875
$ws.Range("M73").Value = -2939
$ws.Range("H94").Value = 103994.5
$ws.Range("J94").Value = 103994.5
$ws.Range("L94").Value = 103994.5
$ws.Range("N94").Value = -105346.5
$ws.Range("H99").Value = 62499.5
$ws.Range("I99").Value = 20000
$ws.Range("K99").Value = 20000
$ws.Range("M99").Value = -17754
$ws.Range("H126").Value = 4152.7085
$ws.Range("I126").Value = 4046.6667
$ws.Range("J126").Value = 4258.75
$ws.Range("K126").Value = 12140.0001
$ws.Range("L126").Value = 12776.25
$ws.Range("M126").Value = -9670.000100000001
$ws.Range("N126").Value = -17716.25
$ws.Range("H132").Value = 3055.5789
$ws.Range("I132").Value = 3055.5789
$ws.Range("K132").Value = 9166.736699999999
$ws.Range("M132").Value = -6636.736699999999
$ws.Range("H139").Value = 86300
$ws.Range("J139").Value = 86300
$ws.Range("L139").Value = 86300
$ws.Range("N139").Value = -96580

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2635.7273
$ws.Range("I16").Value = 2499.3
$ws.Range("J16").Value = 4000
$ws.Range("K16").Value = 2499.3
$ws.Range("L16").Value = 4000
$ws.Range("M16").Value = -2329.3
$ws.Range("N16").Value = -4340
$ws.Range("H22").Value = 5000
$ws.Range("I22").Value = 5000
$ws.Range("K22").Value = 5000
$ws.Range("M22").Value = -4705
$ws.Range("H27").Value = 5000
$ws.Range("I27").Value = 5000
$ws.Range("K27").Value = 5000
$ws.Range("M27").Value = -4893
$ws.Range("H55").Value = 330.35135
$ws.Range("I55").Value = 385.6
$ws.Range("J55").Value = 292.68182
$ws.Range("K55").Value = 385.6
$ws.Range("L55").Value = 292.68182
$ws.Range("M55").Value = -212.6
$ws.Range("N55").Value = -638.68182
$ws.Range("H122").Value = 20282.889
$ws.Range("I122").Value = 20282.889
$ws.Range("K122").Value = 60848.667
$ws.Range("M122").Value = -58398.667
$ws.Range("H132").Value = 502968.06
$ws.Range("I132").Value = 591198
$ws.Range("K132").Value = 1773594
$ws.Range("M132").Value = -1771064
$ws.Range("H136").Value = 4625.75
$ws.Range("I136").Value = 4169.3335
$ws.Range("J136").Value = 4899.6
$ws.Range("K136").Value = 12508.0005
$ws.Range("L136").Value = 14698.8
$ws.Range("M136").Value = -9958.000499999998
$ws.Range("N136").Value = -19798.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 30085
$ws.Range("I58").Value = 30085
$ws.Range("K58").Value = 30085
$ws.Range("M58").Value = -29777
$ws.Range("H132").Value = 38926.758
$ws.Range("I132").Value = 45661.293
$ws.Range("J132").Value = 6601
$ws.Range("K132").Value = 136983.879
$ws.Range("L132").Value = 19803
$ws.Range("M132").Value = -134453.879
$ws.Range("N132").Value = -24863
$ws.Range("H136").Value = 25938.191
$ws.Range("I136").Value = 1582.7878
$ws.Range("K136").Value = 4748.3634
$ws.Range("M136").Value = -2198.3634
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
